$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 430490.44
$ws.Range("J17").Value = 483449.4
$ws.Range("L17").Value = 1450348.2
$ws.Range("N17").Value = -1450684.2

$ws.Range("H51").Value = 2692
$ws.Range("I51").Value = 2575.5
$ws.Range("J51").Value = 2709.2593
$ws.Range("K51").Value = 2575.5
$ws.Range("L51").Value = 2709.2593
$ws.Range("M51").Value = -2091.5
$ws.Range("N51").Value = -3677.2593

$ws.Range("H129").Value = 797.1
$ws.Range("I129").Value = 570
$ws.Range("J129").Value = 919.38464
$ws.Range("K129").Value = 1710
$ws.Range("L129").Value = 2758.15392
$ws.Range("M129").Value = 3290
$ws.Range("N129").Value = -12758.15392

$ws.Range("H135").Value = 344.125
$ws.Range("I135").Value = 309.61404
$ws.Range("J135").Value = 625.1429000000001
$ws.Range("K135").Value = 2786.52636
$ws.Range("L135").Value = 5626.2861
$ws.Range("M135").Value = -251.5263599999998
$ws.Range("N135").Value = -10696.2861

$ws.Range("H137").Value = 2851.5342
$ws.Range("I137").Value = 1133.091
$ws.Range("J137").Value = 5458.8276
$ws.Range("K137").Value = 3399.273
$ws.Range("L137").Value = 16376.4828
$ws.Range("M137").Value = -849.2729999999997
$ws.Range("N137").Value = -21476.4828

$ws.Range("H138").Value = 2109
$ws.Range("I138").Value = 967.5833
$ws.Range("J138").Value = 3935.2666
$ws.Range("K138").Value = 2902.7499
$ws.Range("L138").Value = 11805.7998
$ws.Range("M138").Value = 2237.2501
$ws.Range("N138").Value = -22085.7998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10105735
$ws.Range("I32").Value = 10992920
$ws.Range("J32").Value = 14003.5
$ws.Range("K32").Value = 10992920
$ws.Range("L32").Value = 14003.5
$ws.Range("M32").Value = -10992633
$ws.Range("N32").Value = -14577.5

$ws.Range("H61").Value = 1530.5358
$ws.Range("I61").Value = 1369.3478
$ws.Range("J61").Value = 2272
$ws.Range("K61").Value = 1369.3478
$ws.Range("L61").Value = 2272
$ws.Range("M61").Value = -1157.3478
$ws.Range("N61").Value = -2696

$ws.Range("H74").Value = 5693.8096
$ws.Range("I74").Value = 826.82355
$ws.Range("K74").Value = 826.82355
$ws.Range("M74").Value = 47.17645000000005

$ws.Range("H77").Value = 5693.8096
$ws.Range("I77").Value = 826.82355
$ws.Range("K77").Value = 4134.117749999999
$ws.Range("M77").Value = 233.8822500000006

$ws.Range("H122").Value = 820.05884
$ws.Range("I122").Value = 681.7
$ws.Range("J122").Value = 1017.7143
$ws.Range("K122").Value = 2045.1
$ws.Range("L122").Value = 3053.1429
$ws.Range("M122").Value = 404.8999999999999
$ws.Range("N122").Value = -7953.1429

$ws.Range("H136").Value = 1530.5358
$ws.Range("I136").Value = 1369.3478
$ws.Range("J136").Value = 2272
$ws.Range("K136").Value = 4108.0434
$ws.Range("L136").Value = 6816
$ws.Range("M136").Value = -1558.0434
$ws.Range("N136").Value = -11916

$ws.Range("H139").Value = 133483.38
$ws.Range("I139").Value = 47590
$ws.Range("J139").Value = 145753.86
$ws.Range("K139").Value = 47590
$ws.Range("L139").Value = 145753.86
$ws.Range("M139").Value = -42450
$ws.Range("N139").Value = -156033.86

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 14490.167
$ws.Range("I26").Value = 9605.5
$ws.Range("J26").Value = 24259.5
$ws.Range("K26").Value = 9605.5
$ws.Range("L26").Value = 24259.5
$ws.Range("M26").Value = -9313.5
$ws.Range("N26").Value = -24843.5

$ws.Range("H80").Value = 413.70587
$ws.Range("I80").Value = 337.55554
$ws.Range("J80").Value = 499.375
$ws.Range("K80").Value = 337.55554
$ws.Range("L80").Value = 499.375
$ws.Range("M80").Value = 660.4444599999999
$ws.Range("N80").Value = -2495.375

$ws.Range("H81").Value = 20052
$ws.Range("J81").Value = 20052
$ws.Range("L81").Value = 20052
$ws.Range("N81").Value = -22174

$ws.Range("H83").Value = 413.70587
$ws.Range("I83").Value = 337.55554
$ws.Range("J83").Value = 499.375
$ws.Range("K83").Value = 1687.7777
$ws.Range("L83").Value = 2496.875
$ws.Range("M83").Value = 3304.2223
$ws.Range("N83").Value = -12480.875

$ws.Range("H84").Value = 20052
$ws.Range("J84").Value = 20052
$ws.Range("L84").Value = 60156
$ws.Range("N84").Value = -70764

$ws.Range("H134").Value = 1305.125
$ws.Range("I134").Value = 1121.1818
$ws.Range("J134").Value = 1709.8
$ws.Range("K134").Value = 3363.5454
$ws.Range("L134").Value = 5129.4
$ws.Range("M134").Value = -828.5454
$ws.Range("N134").Value = -10199.4

$ws.Range("H135").Value = 35390
$ws.Range("J135").Value = 35390
$ws.Range("L135").Value = 35390
$ws.Range("N135").Value = -45530

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15659987
$ws.Range("I31").Value = 26316972
$ws.Range("J31").Value = 84392.53999999999
$ws.Range("K31").Value = 26316972
$ws.Range("L31").Value = 84392.53999999999
$ws.Range("M31").Value = -26316677
$ws.Range("N31").Value = -84982.53999999999

$ws.Range("H34").Value = 15659987
$ws.Range("I34").Value = 26316972
$ws.Range("J34").Value = 84392.53999999999
$ws.Range("K34").Value = 26316972
$ws.Range("L34").Value = 84392.53999999999
$ws.Range("M34").Value = -26316770
$ws.Range("N34").Value = -84796.53999999999

$ws.Range("H58").Value = 1480.6765
$ws.Range("I58").Value = 1346.2307
$ws.Range("J58").Value = 1563.9048
$ws.Range("K58").Value = 1346.2307
$ws.Range("L58").Value = 1563.9048
$ws.Range("M58").Value = -1143.2307
$ws.Range("N58").Value = -1969.9048

$ws.Range("H105").Value = 2061.261
$ws.Range("I105").Value = 1289.2667
$ws.Range("J105").Value = 3508.75
$ws.Range("K105").Value = 1289.2667
$ws.Range("L105").Value = 3508.75
$ws.Range("M105").Value = 457.7333000000001
$ws.Range("N105").Value = -7002.75

$ws.Range("H122").Value = 1111.091
$ws.Range("I122").Value = 1023.1111
$ws.Range("J122").Value = 1507
$ws.Range("K122").Value = 3069.3333
$ws.Range("L122").Value = 4521
$ws.Range("M122").Value = -619.3332999999998
$ws.Range("N122").Value = -9421

$ws.Range("H134").Value = 1891.75
$ws.Range("I134").Value = 1881.6666
$ws.Range("J134").Value = 1922
$ws.Range("K134").Value = 5644.9998
$ws.Range("L134").Value = 5766
$ws.Range("M134").Value = -3109.9998
$ws.Range("N134").Value = -10836

$ws.Range("H136").Value = 1480.6765
$ws.Range("I136").Value = 1346.2307
$ws.Range("J136").Value = 1563.9048
$ws.Range("K136").Value = 4038.6921
$ws.Range("L136").Value = 4691.7144
$ws.Range("M136").Value = -1488.6921
$ws.Range("N136").Value = -9791.714400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 483.51724
$ws.Range("I5").Value = 253.94118
$ws.Range("K5").Value = 761.82354
$ws.Range("M5").Value = -649.82354

$ws.Range("H12").Value = 39.40625
$ws.Range("I12").Value = 2.2857144
$ws.Range("J12").Value = 49.8
$ws.Range("K12").Value = 6.857143199999999
$ws.Range("L12").Value = 149.4
$ws.Range("M12").Value = 166.1428568
$ws.Range("N12").Value = -495.4

$ws.Range("H88").Value = 4785.7144
$ws.Range("J88").Value = 4785.7144
$ws.Range("L88").Value = 14357.1432
$ws.Range("N88").Value = -15213.1432

$ws.Range("H91").Value = 4785.7144
$ws.Range("J91").Value = 4785.7144
$ws.Range("L91").Value = 14357.1432
$ws.Range("N91").Value = -17321.1432

$ws.Range("H122").Value = 1041.0889
$ws.Range("I122").Value = 499.16666
$ws.Range("J122").Value = 1402.3704
$ws.Range("K122").Value = 4492.49994
$ws.Range("L122").Value = 12621.3336
$ws.Range("M122").Value = -2042.49994
$ws.Range("N122").Value = -17521.3336

$ws.Range("H132").Value = 4348534
$ws.Range("I132").Value = 12500452
$ws.Range("J132").Value = 844.2
$ws.Range("K132").Value = 112504068
$ws.Range("L132").Value = 7597.8
$ws.Range("M132").Value = -112501538
$ws.Range("N132").Value = -12657.8

$ws.Range("H135").Value = 483.51724
$ws.Range("I135").Value = 253.94118
$ws.Range("K135").Value = 2285.47062
$ws.Range("M135").Value = 249.5293799999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1413.8334
$ws.Range("I126").Value = 992
$ws.Range("J126").Value = 1624.75
$ws.Range("K126").Value = 2976
$ws.Range("L126").Value = 4874.25
$ws.Range("M126").Value = -506
$ws.Range("N126").Value = -9814.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2212.5
$ws.Range("I40").Value = 1812
$ws.Range("J40").Value = 2613
$ws.Range("K40").Value = 1812
$ws.Range("L40").Value = 2613
$ws.Range("M40").Value = -1676
$ws.Range("N40").Value = -2885

$ws.Range("H61").Value = 1519.091
$ws.Range("I61").Value = 1315.4546
$ws.Range("J61").Value = 1620.909
$ws.Range("K61").Value = 1315.4546
$ws.Range("L61").Value = 1620.909
$ws.Range("M61").Value = -1113.4546
$ws.Range("N61").Value = -2024.909

$ws.Range("H113").Value = 1519.091
$ws.Range("I113").Value = 1315.4546
$ws.Range("J113").Value = 1620.909
$ws.Range("K113").Value = 1315.4546
$ws.Range("L113").Value = 1620.909
$ws.Range("M113").Value = 854.5454
$ws.Range("N113").Value = -5960.909

$ws.Range("H132").Value = 4263.3076
$ws.Range("I132").Value = 5279.613
$ws.Range("J132").Value = 2763.0476
$ws.Range("K132").Value = 15838.839
$ws.Range("L132").Value = 8289.1428
$ws.Range("M132").Value = -13308.839
$ws.Range("N132").Value = -13349.1428

$ws.Range("H136").Value = 6137.222
$ws.Range("I136").Value = 2354.5454
$ws.Range("J136").Value = 22781
$ws.Range("K136").Value = 7063.6362
$ws.Range("L136").Value = 68343
$ws.Range("M136").Value = -4513.6362
$ws.Range("N136").Value = -73443

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2882.2712
$ws.Range("I136").Value = 3583.5625
$ws.Range("J136").Value = 2051.111
$ws.Range("K136").Value = 10750.6875
$ws.Range("L136").Value = 6153.333
$ws.Range("M136").Value = -8200.6875
$ws.Range("N136").Value = -11253.333
